# "Update countries & provincias Spain"
#
# The upstream data refresh (14-May-2020, now at 22:35) did two things:
#   1. Refreshed the COVID-19 counters for every country.
#   2. Re-sorted a handful of countries that were tied on case count, which
#      shuffled a few neighbouring rows (e.g. "Togo" jumped above "Haiti",
#      "Santo Tome y Principe" & "Madagascar"; "Uganda" jumped above
#      "Guadalupe" ... "Brunei"; "Santa Lucia" and "Belice" swapped).
#
# Below we just write the final, as-published values straight into the
# affected cells (country name in column A plus the Casos
# totales/Nuevos/Activos/Recuperados/Criticos/MuertesHoy/Muertes columns
# B-H for that row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / "last updated" banner
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 22:35"

# Row 4 - Estados Unidos (counts refreshed; country unchanged)
$ws.Cells.Item(4, 2).Value = 1449744
$ws.Cells.Item(4, 3).Value = 19396
$ws.Cells.Item(4, 4).Value = 316169
$ws.Cells.Item(4, 5).Value = 1046995
$ws.Cells.Item(4, 7).Value = 1383
$ws.Cells.Item(4, 8).Value = 86580

# Rows 143-146: "Togo" moves up above Haiti / Santo Tome y Principe / Madagascar
# Row 143 - was Haiti, now Togo
$ws.Cells.Item(143, 1).Value = "Togo"
$ws.Cells.Item(143, 2).Value = 238
$ws.Cells.Item(143, 3).Value = 19
$ws.Cells.Item(143, 4).Value = 96
$ws.Cells.Item(143, 5).Value = 131
$ws.Cells.Item(143, 8).Value = 11

# Row 144 - was Madagascar, now Haiti
$ws.Cells.Item(144, 1).Value = "Haiti"
$ws.Cells.Item(144, 2).Value = 234
$ws.Cells.Item(144, 3).Value = 15
$ws.Cells.Item(144, 4).Value = 17
$ws.Cells.Item(144, 5).Value = 199
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 8).Value = 18

# Row 145 - was Santo Tome y Principe, now Madagascar
$ws.Cells.Item(145, 1).Value = "Madagascar"
$ws.Cells.Item(145, 2).Value = 230
$ws.Cells.Item(145, 3).Value = 18
$ws.Cells.Item(145, 4).Value = 108
$ws.Cells.Item(145, 5).Value = 122
$ws.Cells.Item(145, 6).Value = 1
$ws.Cells.Item(145, 8).Value = 0

# Row 146 - was Togo, now Santo Tome y Principe
$ws.Cells.Item(146, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(146, 2).Value = 220
$ws.Cells.Item(146, 4).Value = 4
$ws.Cells.Item(146, 5).Value = 210
$ws.Cells.Item(146, 8).Value = 6

# Rows 153-158: "Uganda" moves up above Guadalupe / Guayana Francesa /
# Gibraltar / Republica de Africa Central / Brunei
# Row 153 - was Guadalupe, now Uganda
$ws.Cells.Item(153, 1).Value = "Uganda"
$ws.Cells.Item(153, 2).Value = 160
$ws.Cells.Item(153, 3).Value = 21
$ws.Cells.Item(153, 4).Value = 63
$ws.Cells.Item(153, 5).Value = 97
$ws.Cells.Item(153, 6).Value = 0
$ws.Cells.Item(153, 8).Value = 0

# Row 154 - was Guayana Francesa, now Guadalupe
$ws.Cells.Item(154, 1).Value = "Guadalupe"
$ws.Cells.Item(154, 2).Value = 155
$ws.Cells.Item(154, 4).Value = 109
$ws.Cells.Item(154, 5).Value = 33
$ws.Cells.Item(154, 6).Value = 4
$ws.Cells.Item(154, 8).Value = 13

# Row 155 - was Gibraltar, now Guayana Francesa
$ws.Cells.Item(155, 1).Value = "Guayana Francesa"
$ws.Cells.Item(155, 2).Value = 153
$ws.Cells.Item(155, 4).Value = 124
$ws.Cells.Item(155, 5).Value = 28
$ws.Cells.Item(155, 8).Value = 1

# Row 156 - was Republica de Africa Central, now Gibraltar
$ws.Cells.Item(156, 1).Value = "Gibraltar"
$ws.Cells.Item(156, 2).Value = 147
$ws.Cells.Item(156, 4).Value = 144
$ws.Cells.Item(156, 5).Value = 3

# Row 157 - was Brunei, now Republica de Africa Central
$ws.Cells.Item(157, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(157, 2).Value = 143
$ws.Cells.Item(157, 4).Value = 10
$ws.Cells.Item(157, 5).Value = 133
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 8).Value = 0

# Row 158 - was Uganda, now Brunei
$ws.Cells.Item(158, 1).Value = "Brunei"
$ws.Cells.Item(158, 2).Value = 141
$ws.Cells.Item(158, 4).Value = 134
$ws.Cells.Item(158, 5).Value = 6
$ws.Cells.Item(158, 6).Value = 2
$ws.Cells.Item(158, 8).Value = 1

# Rows 194-195: "Santa Lucia" and "Belice" swap places
# Row 194 - was Belice, now Santa Lucia
$ws.Cells.Item(194, 1).Value = "Santa Lucia"
$ws.Cells.Item(194, 4).Value = 18
$ws.Cells.Item(194, 8).Value = 0

# Row 195 - was Santa Lucia, now Belice
$ws.Cells.Item(195, 1).Value = "Belice"
$ws.Cells.Item(195, 4).Value = 16
$ws.Cells.Item(195, 8).Value = 2
